# Update weekly excess mortality analysis ("Berekening oversterfte CBS.xlsx")
# - revises a handful of "Waargenomen" (observed) weekly mortality figures
#   in column G (which ripple into the "Oversterfte" column I via =G-H)
# - appends a new week (week 34) as row 26
# - extends the week 11-33/34 summary totals in row 28 to include the new row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Revised observed weekly mortality figures (column G) ---
# Column I (Oversterfte = G - H) recalculates automatically.
$ws.Range("G11").Value = 2982
$ws.Range("G12").Value = 2773
$ws.Range("G13").Value = 2769
$ws.Range("G17").Value = 2691
$ws.Range("G19").Value = 2635
$ws.Range("G20").Value = 2612
$ws.Range("G21").Value = 2520
$ws.Range("G22").Value = 2668
$ws.Range("G23").Value = 2641
$ws.Range("G24").Value = 2625
$ws.Range("G25").Value = 3189

# --- New week 34 row ---
$ws.Range("F26").Value = 34
$ws.Range("G26").Value = 2810
$ws.Range("H26").Value = 3022
$ws.Range("I26").Formula = "=G26-H26"

# --- Extend the totals row to cover the new row ---
$ws.Range("G28").Formula = "=SUM(G3:G26)"
$ws.Range("H28").Formula = "=SUM(H3:H26)"
$ws.Range("I28").Formula = "=SUM(I3:I26)"

# --- Restore the author's last active-cell selection ---
$ws.Range("K20").Select()

$wb.Save()
